{"js": "// The \"Shen, W., Darrel Jenerette...\" reference paragraph currently sits\n// immediately before the paragraph that only contains the _GoBack\n// bookmark. The edit moves that reference text so it lives *inside* the\n// bookmark paragraph (appended right after the bookmark), empties out the\n// paragraph that used to hold it, and adds a new blank paragraph right\n// after the (now populated) bookmark paragraph.\n\nconst bookmarkRange = context.document.getBookmarkRange(\"_GoBack\");\nconst bookmarkPara = bookmarkRange.paragraphs.getFirst();\nconst citationPara = bookmarkPara.getPrevious();\n\ncitationPara.load(\"text\");\nawait context.sync();\n\n// Split the citation text the same way the original run layout did, so\n// the destination paragraph ends up with the same run boundaries.\nconst segments = [\n  \"Shen\",\n  \", W., Darrel \",\n  \"Jenerette\",\n  \", G., Wu, J. and H Gardner, R., 2004.\",\n  \" Evaluating empirical scaling relations of pattern metrics with simulated landscapes. \",\n  \"Ecography\",\n  \", 27(4), pp.459-469.\"\n];\n\nfor (const segment of segments) {\n  bookmarkPara.insertText(segment, Word.InsertLocation.end);\n}\n\n// Empty out the old paragraph - its text has moved into bookmarkPara.\ncitationPara.clear();\n\n// Add the new blank paragraph right after the bookmark paragraph.\nbookmarkPara.insertParagraph(\"\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# The \"Shen, W., Darrel Jenerette...\" reference paragraph currently sits\n# immediately before the paragraph that only contains the _GoBack\n# bookmark. This edit moves that reference text so it lives *inside* the\n# bookmark paragraph (appended right after the bookmark), empties out the\n# paragraph that used to hold it, and adds a new blank paragraph right\n# after the (now populated) bookmark paragraph.\n\n$d = $word.ActiveDocument\n\n$bm = $d.Bookmarks.Item(\"_GoBack\")\n$bmPara = $bm.Range.Paragraphs.Item(1)\n$bmIndex = $bmPara.Index\n$citationPara = $d.Paragraphs.Item($bmIndex - 1)\n\n# Clear the citation paragraph's text. Building the Range explicitly via\n# $d.Range(start, end) (rather than reusing the Paragraph's .Range\n# property) is what makes the whole multi-run paragraph clear in one\n# shot here.\n$citationStart = $citationPara.Range.Start\n$citationEnd = $citationPara.Range.End\n$clearRange = $d.Range($citationStart, $citationEnd)\n$clearRange.Text = \"\"\n\n# Re-fetch the bookmark paragraph - text positions shifted after the\n# clear above.\n$bm = $d.Bookmarks.Item(\"_GoBack\")\n$bmPara = $bm.Range.Paragraphs.Item(1)\n\n# Re-insert the citation text into the bookmark paragraph, split into the\n# same run boundaries the original paragraph used, tracking the insertion\n# point explicitly so each segment lands after the previous one.\n$segments = @(\n  \"Shen\",\n  \", W., Darrel \",\n  \"Jenerette\",\n  \", G., Wu, J. and H Gardner, R., 2004.\",\n  \" Evaluating empirical scaling relations of pattern metrics with simulated landscapes. \",\n  \"Ecography\",\n  \", 27(4), pp.459-469.\"\n)\n\n$insertPoint = $bmPara.Range.End\nforeach ($segment in $segments) {\n  $insertRange = $d.Range($insertPoint, $insertPoint)\n  $insertRange.InsertAfter($segment)\n  $insertPoint = $insertRange.End\n}\n\n# Add the new blank paragraph right after the bookmark paragraph.\n$bm = $d.Bookmarks.Item(\"_GoBack\")\n$bmPara = $bm.Range.Paragraphs.Item(1)\n$bmPara.Range.InsertParagraphAfter()\n"}
